# Update the "想去人数" (interested-count) column F on both the "展览"
# and "全部类型" sheets with refreshed values.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    3  = 508
    4  = 439
    5  = 8527
    7  = 1513
    8  = 176
    11 = 252
    12 = 390
    13 = 247
    15 = 31
    16 = 129
    18 = 459
    19 = 1235
    20 = 197
    21 = 78
    22 = 134
    23 = 94
    24 = 124
    25 = 74
    26 = 110
    27 = 101
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
